# "small preview of the website + stickers added"
#
# The underlying sheet data only has two real content changes:
#   B2: the Miyabi/Yixuan/Alice stickers line used to be one run separated
#       by "; " -- it is now wrapped onto separate lines (one name per line),
#       matching how every other multi-name cell in the sheet is formatted.
#   C3: "Yanagi(Anomaly, Disorder);" loses the ", Disorder" qualifier and
#       becomes simply "Yanagi(Anomaly);".
#
# (Everything else in the source diff -- the shared-string table being
# reshuffled, cell <v> indices changing -- is just a side effect of Excel
# rewriting xl/sharedStrings.xml after these edits; the actual cell
# contents for every other cell are unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Miyabi(Anomaly, Quickswap);`nYixuan(Sheer, Quickswap);`nAlice(Anomaly, Quickswap);[+1]"
$ws.Range("C3").Value = "Yanagi(Anomaly);"

# Cosmetic: the author also scrolled the sheet over and re-selected C3
# (a small preview pane of the site was open alongside Excel).
$ws.Range("C3").Select()
